# -----------------------------------------------------------------------------
# Updates the cryptocurrency price/volume snapshot on Sheet1 (cryptos.xlsx),
# matching the "Updated cryptos list ... with GitHub Actions" commit.
#
# For each affected row, the Price (column D) and/or Volume(1h) (column E)
# text is refreshed with freshly scraped values. A handful of rows
# (35-38 and 51) also get new Coin name / Link text because the coin
# ranking reshuffled, bringing in a different coin (e.g. Aave -> Tezos).
#
# All of these cells store plain TEXT in the workbook (e.g. "1.001",
# "  -0.04%  "), not numbers. Excel's COM layer auto-converts a plain
# numeric-looking string assigned to .Value into a real number (losing the
# trailing zeros / introducing floating point noise), so for any
# replacement value that parses as a number we first force the cell's
# NumberFormat to Text ("@"), assign the value, and then reset the cell
# style back to "Normal" so the cell keeps its original (unstyled) look
# instead of being left tagged with the temporary Text format.
# -----------------------------------------------------------------------------

function Set-TextValue($ws, $cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.566.22'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.751.49'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("E4").Value = '  +0.11%  '
Set-TextValue $ws 'D5' '324.10'
$ws.Range("E5").Value = '  -0.18%  '
Set-TextValue $ws 'D6' '1.001'
$ws.Range("E6").Value = '  +0.06%  '
Set-TextValue $ws 'D7' '0.4491'
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("E8").Value = '  -1.97%  '
Set-TextValue $ws 'D9' '0.07455'
$ws.Range("E9").Value = '  -0.93%  '
Set-TextValue $ws 'D10' '41.35'
$ws.Range("E10").Value = '  -1.74%  '
Set-TextValue $ws 'D11' '1.081'
$ws.Range("E11").Value = '  -2.40%  '
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("E13").Value = '  -0.43%  '
Set-TextValue $ws 'D14' '5.975'
$ws.Range("E14").Value = '  -1.46%  '
Set-TextValue $ws 'D15' '7.145'
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("D16").Value = '1.744.45'
$ws.Range("E16").Value = '  -0.75%  '
Set-TextValue $ws 'D17' '93.69'
$ws.Range("E17").Value = '  +0.80%  '
$ws.Range("E18").Value = '  -0.97%  '
Set-TextValue $ws 'D19' '0.06364'
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("E20").Value = '  +0.01%  '
Set-TextValue $ws 'D21' '17.11'
$ws.Range("E21").Value = '  -0.22%  '
Set-TextValue $ws 'D22' '5.731'
$ws.Range("E22").Value = '  -2.45%  '
$ws.Range("D23").Value = '27.615.62'
$ws.Range("E23").Value = '  -0.03%  '
Set-TextValue $ws 'D24' '11.18'
$ws.Range("E24").Value = '  -0.79%  '
Set-TextValue $ws 'D25' '2.083'
$ws.Range("E25").Value = '  -0.19%  '
Set-TextValue $ws 'D26' '165.44'
$ws.Range("E26").Value = '  +1.42%  '
Set-TextValue $ws 'D27' '20.14'
$ws.Range("E27").Value = '  -1.76%  '
$ws.Range("D28").Value = '1.949.90'
$ws.Range("E28").Value = '  -0.36%  '
Set-TextValue $ws 'D29' '2.100'
$ws.Range("E29").Value = '  -1.90%  '
Set-TextValue $ws 'D30' '125.21'
$ws.Range("E30").Value = '  -0.60%  '
Set-TextValue $ws 'D31' '1.083'
$ws.Range("E31").Value = '  -1.49%  '
Set-TextValue $ws 'D32' '0.09175'
$ws.Range("E32").Value = '  +1.12%  '
Set-TextValue $ws 'D33' '3.650'
$ws.Range("E33").Value = '  +0.46%  '
Set-TextValue $ws 'D34' '5.497'
$ws.Range("E34").Value = '  -1.31%  '
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws 'D35' '0.02285'
$ws.Range("E35").Value = '  -1.04%  '
$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws 'D36' '11.73'
$ws.Range("E36").Value = '  -3.82%  '
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws 'D37' '0.2086'
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws 'D38' '0.06006'
$ws.Range("E38").Value = '  -0.23%  '
Set-TextValue $ws 'D39' '0.6269'
$ws.Range("E39").Value = '  -2.56%  '
Set-TextValue $ws 'D40' '4.919'
$ws.Range("E40").Value = '  -0.59%  '
Set-TextValue $ws 'D41' '1.179'
$ws.Range("E41").Value = '  -1.19%  '
Set-TextValue $ws 'D42' '1.397'
$ws.Range("E42").Value = '  -0.07%  '
Set-TextValue $ws 'D43' '7.766'
$ws.Range("E43").Value = '  -1.34%  '
Set-TextValue $ws 'D44' '13.14'
$ws.Range("E44").Value = '  -0.98%  '
Set-TextValue $ws 'D45' '3.715'
$ws.Range("E45").Value = '  -0.04%  '
Set-TextValue $ws 'D46' '0.5854'
$ws.Range("E46").Value = '  -1.31%  '
Set-TextValue $ws 'D47' '121.74'
$ws.Range("E47").Value = '  -0.43%  '
$ws.Range("E48").Value = '  -2.50%  '
Set-TextValue $ws 'D49' '0.06879'
$ws.Range("E49").Value = '  +0.08%  '
Set-TextValue $ws 'D50' '1.127'
$ws.Range("E50").Value = '  -3.55%  '
$ws.Range("B51").Value = 'Tezos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz'
Set-TextValue $ws 'D51' '1.125'
$ws.Range("E51").Value = '  -0.55%  '
